# Append a new block of paragraphs (text paragraph followed by a blank
# paragraph, repeated) to the end of the document, right after the last
# existing (empty) paragraph and before the section break.

$d = $word.ActiveDocument

$paragraphs = @(
    "Encontre maneiras de apoiar o relacionamento entre jogador e personagem e construir confiança e intimidade ao longo do tempo.",
    "",
    "Se você estiver criando um mundo ficcional repleto de histórias, certifique-se de que haja oportunidades para “tensão e liberação”.",
    "",
    "Considere a agência do seu jogador e, se houver falta de agência, torne isso significativo também.",
    "",
    "Ofereça oportunidades de reflexão e vínculo.",
    "",
    "Considere como o ponto de vista (primeira pessoa versus terceira pessoa) pode afetar a empatia pelos jogadores e personagens.",
    "",
    "Avalie o papel da narrativa e do envolvimento narrativo (ou “transporte”) nos jogos e sua relação com o sentimento de “agência” dos jogadores, a tomada de perspectiva e a construção de relacionamentos.",
    "",
    "Considere e avalie como diferentes contextos, públicos e experiências e expectativas anteriores influenciam a empatia por meio dos jogos."
)

foreach ($text in $paragraphs) {
    $lastPara = $d.Paragraphs.Last
    $lastPara.Range.InsertParagraphAfter()
    if ($text -ne "") {
        $d2 = $word.ActiveDocument
        $newPara = $d2.Paragraphs.Last
        $newPara.Range.InsertBefore($text)
    }
}
